$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings must be introduced in this exact order so the
#     rebuilt sharedStrings.xml table lands on the same indices as the
#     target workbook (index is assigned in first-use order). ---
$ws.Range("J4").Value2 = "val acc after 5 epochs"
$ws.Range("A3").Value2 = "Structured Pruning After Training"
$ws.Range("G3").Value2 = "Structured Pruning with Training"
$ws.Range("J17").Value2 = "97.91% -> 11.35%"
$ws.Range("J16").Value2 = "97.83% -> 95.57%"
$ws.Range("J15").Value2 = "97.87% -> 96.93%"
$ws.Range("J14").Value2 = "97.73% -> 97.69%"

# --- "cpu" block (rows 5-10): updated post_acc numbers ---
$ws.Range("J5").Value2 = 0.97970000000000002
$ws.Range("J6").Value2 = 0.98060000000000003
$ws.Range("J7").Value2 = 0.97829999999999995
$ws.Range("J8").Value2 = 0.97950000000000004
$ws.Range("J9").Value2 = 0.98199999999999998
$ws.Range("J10").Value2 = 0.97740000000000005

# --- rows 12-13: J values cleared (kept style, no content) ---
$ws.Range("J12").ClearContents()
$ws.Range("J13").ClearContents()

# --- rows 19-20: J values cleared ---
$ws.Range("J19").ClearContents()
$ws.Range("J20").ClearContents()

# --- rows 21-22: updated J values ---
$ws.Range("J21").Value2 = 0.98040000000000005
$ws.Range("J22").Value2 = 0.97940000000000005

# --- rows 23-24: J values cleared ---
$ws.Range("J23").ClearContents()
$ws.Range("J24").ClearContents()

# --- rows 26-31: J values cleared ---
$ws.Range("J26").ClearContents()
$ws.Range("J27").ClearContents()
$ws.Range("J28").ClearContents()
$ws.Range("J29").ClearContents()
$ws.Range("J30").ClearContents()
$ws.Range("J31").ClearContents()

# --- column widths: widen column I and give column J an explicit width ---
$ws.Columns.Item(9).ColumnWidth = 12.333333333333334
$ws.Columns.Item(10).ColumnWidth = 18

# --- view: scroll/selection moved to J19 ---
$ws.Range("J19").Select()
